function Split-RunText {
    param($doc, $needle, $segs)
    $rng = $doc.Content
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $needle"
        return
    }
    $start = $rng.Start
    $newText = ""
    foreach ($s in $segs) { $newText = $newText + $s }
    $rng.Text = $newText
    $pos = $start
    for ($i = 0; $i -lt $segs.Count - 1; $i++) {
        $segLen = $segs[$i].Length
        $segStart = $pos
        $segEnd = $pos + $segLen
        $sub = $doc.Range($segStart, $segEnd)
        $sub.Font.Bold = $true
        $sub.Font.Bold = $false
        $pos = $segEnd
    }
}

$d = $word.ActiveDocument

Split-RunText $d "2.a) " @("2", ".", "a ")
Split-RunText $d "2.b) SaaS" @("2.b", " ", "SaaS")
Split-RunText $d "2.c) No tiene acceso a los Servicios de la Compañía" @("2.c", " ", " No tiene acceso a los Servicios de la Compañía")
Split-RunText $d "3.a) SI" @("3.a", " ", " SI")
Split-RunText $d "3.b) NO" @("3.b", " ", " ", "NO")
Split-RunText $d "4.a) SI" @("4.a", " ", " ", " SI")
Split-RunText $d "4.b) NO" @("4.b", " ", " ", " NO")
Split-RunText $d "5.a) " @("5.a", " ", " ")
Split-RunText $d "5.b) En la Nube" @("5.b", " ", " En la Nube")
Split-RunText $d "5.c) Respaldo en Almacenamiento Físicos Fuera de las premisas" @("5.c", " ", " ", " Respaldo en Almacenamiento Físicos Fuera de las premisas")
Split-RunText $d "6.a) VPN" @("6.a", " ", " ", " VPN")
Split-RunText $d "6.b) Correo Electrónico" @("6.b", " ", " Correo Electrónico")
Split-RunText $d "6.c) ERP" @("6.c", " ", " ERP")
Split-RunText $d "6.d) Autorizaciones Financieras" @("6.d", " ", " Autorizaciones Financieras")
Split-RunText $d "6.e) Lo desconozco" @("6.e", " ", " Lo desconozco")
Split-RunText $d "7.a) SI" @("7.a", " ", "SI")
Split-RunText $d "7.b) NO" @("7.b", " ", "NO")
Split-RunText $d "8.a) SI" @("8.a", " ", " SI")
Split-RunText $d "8.b) NO" @("8.b", " ", "NO")
Split-RunText $d "9.a) AVS" @("9.a AVS")
Split-RunText $d "9.b) IPS" @("9.b", " ", " IPS")
Split-RunText $d "9.c) URL " @("9.c", " ", " URL ")
Split-RunText $d "9.d) " @("9.d", " ")
Split-RunText $d "9.e) " @("9.e", " ", " ")
Split-RunText $d "9.f) Lo desconozco" @("9.f Lo desconozco")
Split-RunText $d "10.a) SI" @("10.a", " ", " SI")
Split-RunText $d "10.b) NO" @("10.b", " ", " NO")
Split-RunText $d "11.a) SI" @("11.a", " ", " SI")
Split-RunText $d "11.b) NO" @("11.b", " ", " NO")
Split-RunText $d "12.a) Manualmente cuando existen Vulnerabilidades de Alto Impacto" @("12.a", " ", " Manualmente cuando existen Vulnerabilidades de Alto Impacto")
Split-RunText $d "12.b) Manualmente después de cada ciclo de Análisis" @("12.b", " ", " Manualmente después de cada ciclo de Análisis")
Split-RunText $d "12.c) Automatizado (Herramientas de Patch Management)" @("12.c", " ", " Automatizado (Herramientas de Patch Management)")
Split-RunText $d "13.a) SI" @("13.a", " ", " SI")
Split-RunText $d "13.b) NO" @("13.b", " ", " NO")
Split-RunText $d "14.a) SI" @("14.a", " ", " SI")
Split-RunText $d "14.b) NO" @("14.b ", "NO")
Split-RunText $d "15.a) SI" @("15.a", " ", " SI")
Split-RunText $d "15.b) NO" @("15.b", " ", " NO")
Split-RunText $d "16.a) NO " @("16.a", " ", " NO ")
Split-RunText $d "16.b) SI. Presupuesto Menor a USD" @("16.b", " ", " SI. Presupuesto Menor a USD")
Split-RunText $d "16.c) SI. Presupuesto Entre USD 100,001 a USD 300,000" @("16.c", " ", " SI. Presupuesto Entre USD 100,001 a USD 300,000")
Split-RunText $d "16.d) SI. Presupuesto Entre USD 300,001 a USD 500,000" @("16.d", " ", " SI. Presupuesto Entre USD 300,001 a USD 500,000")
Split-RunText $d "16.e) SI. Presupuesto Mayor a USD 500,000" @("16.e", " ", " SI. Presupuesto Mayor a USD 500,000")
